$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.002.81"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "1.641.59"
$ws.Range("E3").Value = "  +0.60%  "
$ws.Range("E4").Value = "  +0.40%  "
$ws.Range("D5").Value = "'215.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.0638"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.255"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.40%  "
$ws.Range("D10").Value = "'19.54"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.06%  "
$ws.Range("D11").Value = "'0.0797"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.63%  "
$ws.Range("D12").Value = "1.869.74"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("D14").Value = "1.646.35"
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("E15").Value = "  +0.06%  "
$ws.Range("E16").Value = "  +1.00%  "
$ws.Range("D17").Value = "'63.35"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "26.104.18"
$ws.Range("D20").Value = "'194.08"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.37%  "
$ws.Range("E21").Value = "  -0.55%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("E23").Value = "  -1.07%  "
$ws.Range("D24").Value = "'1.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  +4.63%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").Value = "'6.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.60%  "
$ws.Range("D29").Value = "'15.52"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.78%  "
$ws.Range("E30").Value = "  +0.73%  "
$ws.Range("D31").Value = "'0.0495"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").Value = "  -0.20%  "
$ws.Range("D33").Value = "'3.25"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -0.94%  "
$ws.Range("E35").Value = "  +1.33%  "
$ws.Range("D36").Value = "'0.904"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.49%  "
$ws.Range("D37").Value = "1.130.47"
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("E38").Value = "  -0.97%  "
$ws.Range("E39").Value = "  -0.52%  "
$ws.Range("D40").Value = "'0.0156"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("E41").Value = "  +0.96%  "
$ws.Range("D42").Value = "'99.01"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.21%  "
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("D44").Value = "1.778.59"
$ws.Range("E44").Value = "  +0.72%  "
$ws.Range("E45").Value = "  +4.87%  "
$ws.Range("D46").Value = "'56.51"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.0522"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.24%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").Value = "'1.49"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.50%  "
$ws.Range("D49").Value = "'7.73"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.40%  "
$ws.Range("E50").Value = "  -0.17%  "
$ws.Range("E51").Value = "  +0.30%  "
